$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Row 194 (A="ILLER LOPEZ ROBERTO FERNANDO", B="PAUTA ASTUDILLO JULIO HERNAN")
# is a duplicate of row 193 and gets removed entirely; all rows below shift up by one.
$ws.Rows.Item(194).Delete()

# After the deletion, the former totals row (367) is now row 366.
# Update the totals to reflect the removed row's contribution
# (D decreases by 326.73, G decreases by 1000; other columns unaffected since they were 0).
$ws.Cells.Item(366, 4).Value = 424433.47
$ws.Cells.Item(366, 7).Value = 373790
